# This workbook contains OLS regression summary output (from Python's
# statsmodels, exported via a pipeline) laid out across columns A:Z,
# rows 2:144. Column H holds the "Date:" stamp and column Z holds the
# "Time:" stamp that statsmodels writes into every regression summary.
# The pipeline was re-run, producing a new date/time stamp; this script
# re-applies that re-run to the workbook:
#   - Date: Wed, 13 Apr 2022  ->  Tue, 03 May 2022   (applies to every row)
#   - Time: stamps shift forward by roughly 1h22m, with each row picking
#     up the particular second-granularity timestamp captured when that
#     row's regression finished running.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column H ("Date:") -----------------------------------------------
# Every data row (2-144) shares the same date stamp.
$ws.Range("H2:H144").Value = "Tue, 03 May 2022"

# --- Column Z ("Time:") ------------------------------------------------
# Each data row picks up the new second-level timestamp that corresponds
# to its old one, per the refreshed run.
$ws.Range("Z2").Value      = "08:07:54"
$ws.Range("Z3:Z11").Value  = "08:07:53"
$ws.Range("Z12:Z40").Value = "08:07:54"
$ws.Range("Z41:Z79").Value = "08:07:55"
$ws.Range("Z80").Value     = "08:07:56"
$ws.Range("Z81:Z83").Value = "08:07:55"
$ws.Range("Z84:Z118").Value = "08:07:56"
$ws.Range("Z119").Value    = "08:07:57"
$ws.Range("Z120").Value    = "08:07:56"
$ws.Range("Z121:Z144").Value = "08:07:57"
